$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.12153434753418
$ws.Range("B1").Value = 2.323867797851562
$ws.Range("C1").Value = 2.396320581436157
$ws.Range("D1").Value = 3.003204107284546
$ws.Range("E1").Value = 2.674479007720947
